# "Scaled speed by factor of 10" — update the params/runs inputs and leave
# the UI selection state (active sheet/tab + selected cells) matching where
# the author was last working when they saved.

$wb = $excel.ActiveWorkbook

$wsRuns   = $wb.Worksheets.Item("runs")
$wsParams = $wb.Worksheets.Item("params")

# --- Data edits --------------------------------------------------------

# runs!B2 "Stop Run ID (inclusive)": 15 -> 1
$wsRuns.Range("B2").Value = 1

# params!D2 "Skill Level [0,5]" for run 1: 1 -> 0.8
$wsParams.Range("D2").Value = 0.8

# params!K2 "Port Config (tab name)" for run 1: config1 -> config5
$wsParams.Range("K2").Value = "config5"

# --- Selection / active sheet state ------------------------------------

# "runs" is no longer the active tab; its saved selection moves to B3.
$wsRuns.Range("B3").Select()

# "params" becomes the active tab, selection on K3.
$wsParams.Activate()
$wsParams.Range("K3").Select()
